$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "horas-trabajadas"
$ws.Range("B2").Value = "personas-residentes-viviendas-familiares"
$ws.Range("C2").Value = "ocupacion-1-digito-descripcion"
$ws.Range("D2").Value = "ocupacion-1-digito-codigo"
$ws.Range("E2").Value = "aragon"
